# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the e3fcde4c-... file in both the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Status column for the e3fcde4c-... row (row 3) on every sheet that shows it.
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Error Detail messages explaining the handback/handoff file name mismatch.
$zhcn.Range("P3").Value = "Handback file name: c4blsp14.b3u is different with handoff file name: e3fcde4c-ce14-4706-8443-dcc4918a5a4f.1cd47957e4e946ca91b7d63550c31576c9244d30.zh-cn."
$dede.Range("P3").Value = "Handback file name: c4blsp14.b3u is different with handoff file name: e3fcde4c-ce14-4706-8443-dcc4918a5a4f.1cd47957e4e946ca91b7d63550c31576c9244d30.de-de."

# Widen the Error Detail column (P) so the new messages are readable.
$zhcn.Range("P1").EntireColumn.ColumnWidth = 39.14
$dede.Range("P1").EntireColumn.ColumnWidth = 39.14
